$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (and values) from row 5 into the new row 6, then overwrite values
$ws.Range("A5:D5").Copy($ws.Range("A6"))

$ws.Range("A6").Value = "longTE"
$ws.Range("B6").Value = "SPED KISQLKE KIQQLKQ ENQQLEE ENSQLEY ENQQLEE ENSQLEY"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "c"

$ws.Range("C14").Select() | Out-Null
